$wb = $excel.ActiveWorkbook

# --- mlr_feature_importance ---
$ws = $wb.Worksheets.Item("mlr_feature_importance")
$ws.Range("A2").Value = "t"
$ws.Range("B2").Value = 67755.10993722807
$ws.Range("C2").Value = 6585.627529298183
$ws.Range("A3").Value = "t_COVID19"
$ws.Range("B3").Value = 43261.06665301397
$ws.Range("C3").Value = 5825.464085737522
$ws.Range("A4").Value = "year"
$ws.Range("B4").Value = 15830.6929194741
$ws.Range("C4").Value = 829.3636677799087
$ws.Range("A5").Value = "COVID19"
$ws.Range("B5").Value = 136.2089780628681
$ws.Range("C5").Value = 9.282654805999261
$ws.Range("A6").Value = "AR1"
$ws.Range("B6").Value = 70.39816739466042
$ws.Range("C6").Value = 63.12708823050743
$ws.Range("A7").Value = "AR4"
$ws.Range("B7").Value = 62.60063090212643
$ws.Range("C7").Value = 67.76356090111369
$ws.Range("A8").Value = "AR2"
$ws.Range("B8").Value = 15.43664624076337
$ws.Range("C8").Value = 43.49998691151036
$ws.Range("A9").Value = "monthcos"
$ws.Range("B9").Value = 4.101836276613176
$ws.Range("C9").Value = 12.79691827461691
$ws.Range("A10").Value = "monthsin"
$ws.Range("B10").Value = -2.509808279387653
$ws.Range("C10").Value = 2.326531643258642
$ws.Range("A11").Value = "AR3"
$ws.Range("B11").Value = -61.00055883955211
$ws.Range("C11").Value = 34.75636445408826
$ws.Range("A12").Value = "AR24"
$ws.Range("B12").Value = -75.73340855874122
$ws.Range("C12").Value = 70.31717498184102
$ws.Range("A13").Value = "AR12"
$ws.Range("B13").Value = -77.17576672323048
$ws.Range("C13").Value = 177.7925286633067

# --- knn_feature_importance ---
$ws = $wb.Worksheets.Item("knn_feature_importance")
$ws.Range("A2").Value = "AR3"
$ws.Range("B2").Value = 0.01914788649848504
$ws.Range("C2").Value = 0.002172022450896206
$ws.Range("A3").Value = "AR1"
$ws.Range("B3").Value = 0.01350575788377175
$ws.Range("C3").Value = 0.007379443800092222
$ws.Range("A4").Value = "monthcos"
$ws.Range("B4").Value = 0.01037931409059896
$ws.Range("C4").Value = 0.0008391514453547171
$ws.Range("A5").Value = "AR4"
$ws.Range("B5").Value = 0.007301860277914773
$ws.Range("C5").Value = 0.001483530116066174
$ws.Range("A6").Value = "t"
$ws.Range("B6").Value = 0.001530342038006882
$ws.Range("C6").Value = 0.003339313384419982
$ws.Range("A7").Value = "monthsin"
$ws.Range("B7").Value = 0.001315973501302503
$ws.Range("C7").Value = 0.001477270899923534
$ws.Range("A8").Value = "COVID19"
$ws.Range("B8").Value = 0.000001660967613981157
$ws.Range("C8").Value = 0.000003321935227962314
$ws.Range("A9").Value = "year"
$ws.Range("B9").Value = 0.0
$ws.Range("C9").Value = 0.0
$ws.Range("A10").Value = "t_COVID19"
$ws.Range("B10").Value = -0.001426436648009588
$ws.Range("C10").Value = 0.0005078442193129963
$ws.Range("A11").Value = "AR24"
$ws.Range("B11").Value = -0.00787180877559912
$ws.Range("C11").Value = 0.006397643861755054
$ws.Range("A12").Value = "AR2"
$ws.Range("B12").Value = -0.01480505154217826
$ws.Range("C12").Value = 0.003356773956036957
$ws.Range("A13").Value = "AR12"
$ws.Range("B13").Value = -0.03293518849259014
$ws.Range("C13").Value = 0.004205126138861062

# --- svr_feature_importance ---
$ws = $wb.Worksheets.Item("svr_feature_importance")
$ws.Range("A2").Value = "AR1"
$ws.Range("B2").Value = 57.93337411391035
$ws.Range("C2").Value = 5.960533232545037
$ws.Range("A3").Value = "COVID19"
$ws.Range("B3").Value = -0.0281106148534036
$ws.Range("C3").Value = 0.00183944249541807
$ws.Range("A4").Value = "t"
$ws.Range("B4").Value = -0.2291060705364089
$ws.Range("C4").Value = 0.1822939775697575
$ws.Range("A5").Value = "t_COVID19"
$ws.Range("B5").Value = -0.3707230705266739
$ws.Range("C5").Value = 3.086658110981423
$ws.Range("A6").Value = "year"
$ws.Range("B6").Value = -0.468663340360581
$ws.Range("C6").Value = 0.382130139063335
$ws.Range("A7").Value = "monthcos"
$ws.Range("B7").Value = -0.9585396583160672
$ws.Range("C7").Value = 0.2784834388008868
$ws.Range("A8").Value = "AR2"
$ws.Range("B8").Value = -1.896227521842002
$ws.Range("C8").Value = 0.6219260210649855
$ws.Range("A9").Value = "monthsin"
$ws.Range("B9").Value = -11.71177230716676
$ws.Range("C9").Value = 0.9466176942423992
$ws.Range("A10").Value = "AR3"
$ws.Range("B10").Value = -48.35597222726415
$ws.Range("C10").Value = 5.3492913365567
$ws.Range("A11").Value = "AR4"
$ws.Range("B11").Value = -81.9800038795905
$ws.Range("C11").Value = 3.979854553391653
$ws.Range("A12").Value = "AR12"
$ws.Range("B12").Value = -235.4171784870921
$ws.Range("C12").Value = 9.229763663469296
$ws.Range("A13").Value = "AR24"
$ws.Range("B13").Value = -256.4893070819343
$ws.Range("C13").Value = 11.03493309957789

# --- xgboost_feature_importance ---
$ws = $wb.Worksheets.Item("xgboost_feature_importance")
$ws.Range("A2").Value = "AR1"
$ws.Range("B2").Value = 0.425656659580013
$ws.Range("C2").Value = 0.03460103429184129
$ws.Range("A3").Value = "monthsin"
$ws.Range("B3").Value = 0.4130806191682181
$ws.Range("C3").Value = 0.06547704582888028
$ws.Range("A4").Value = "AR12"
$ws.Range("B4").Value = 0.3120025316588272
$ws.Range("C4").Value = 0.05322608854322328
$ws.Range("A5").Value = "monthcos"
$ws.Range("B5").Value = 0.1622273328557044
$ws.Range("C5").Value = 0.01056345426060904
$ws.Range("A6").Value = "AR4"
$ws.Range("B6").Value = 0.04582453020936974
$ws.Range("C6").Value = 0.01173351828660557
$ws.Range("A7").Value = "AR3"
$ws.Range("B7").Value = 0.02980463461966916
$ws.Range("C7").Value = 0.007218925505183446
$ws.Range("A8").Value = "t_COVID19"
$ws.Range("B8").Value = 0.0007092539524617969
$ws.Range("C8").Value = 0.00164053914382651
$ws.Range("A9").Value = "t"
$ws.Range("B9").Value = 0.0
$ws.Range("C9").Value = 0.0
$ws.Range("A10").Value = "COVID19"
$ws.Range("B10").Value = 0.0
$ws.Range("C10").Value = 0.0
$ws.Range("A11").Value = "year"
$ws.Range("B11").Value = 0.0
$ws.Range("C11").Value = 0.0
$ws.Range("A12").Value = "AR24"
$ws.Range("B12").Value = -0.02205570019047354
$ws.Range("C12").Value = 0.01534093116441737
$ws.Range("A13").Value = "AR2"
$ws.Range("B13").Value = -0.02540253977711861
$ws.Range("C13").Value = 0.01960265785119763

# --- gbt_feature_importance ---
$ws = $wb.Worksheets.Item("gbt_feature_importance")
$ws.Range("A2").Value = "AR12"
$ws.Range("B2").Value = 0.2321775859359667
$ws.Range("C2").Value = 0.06390225421392057
$ws.Range("A3").Value = "AR1"
$ws.Range("B3").Value = 0.2268388115259635
$ws.Range("C3").Value = 0.03153171318251697
$ws.Range("A4").Value = "monthsin"
$ws.Range("B4").Value = 0.1316960835581719
$ws.Range("C4").Value = 0.01952705728908861
$ws.Range("A5").Value = "AR4"
$ws.Range("B5").Value = 0.04732153486811361
$ws.Range("C5").Value = 0.02421835596858743
$ws.Range("A6").Value = "t_COVID19"
$ws.Range("B6").Value = 0.0
$ws.Range("C6").Value = 0.0
$ws.Range("A7").Value = "t"
$ws.Range("B7").Value = 0.0
$ws.Range("C7").Value = 0.0
$ws.Range("A8").Value = "COVID19"
$ws.Range("B8").Value = 0.0
$ws.Range("C8").Value = 0.0
$ws.Range("A9").Value = "year"
$ws.Range("B9").Value = 0.0
$ws.Range("C9").Value = 0.0
$ws.Range("A10").Value = "AR3"
$ws.Range("B10").Value = -0.01916113998238074
$ws.Range("C10").Value = 0.01145586982356775
$ws.Range("A11").Value = "AR2"
$ws.Range("B11").Value = -0.03594260667973947
$ws.Range("C11").Value = 0.01018284653668898
$ws.Range("A12").Value = "monthcos"
$ws.Range("B12").Value = -0.04131414818409249
$ws.Range("C12").Value = 0.01754869719133724
$ws.Range("A13").Value = "AR24"
$ws.Range("B13").Value = -0.06239055970819831
$ws.Range("C13").Value = 0.01464808277542066

# --- elasticnet_feature_importance ---
$ws = $wb.Worksheets.Item("elasticnet_feature_importance")
$ws.Range("A2").Value = "t"
$ws.Range("B2").Value = 0.005456035773996337
$ws.Range("C2").Value = 0.008219407967185982
$ws.Range("A3").Value = "t_COVID19"
$ws.Range("B3").Value = 0.000972830435090355
$ws.Range("C3").Value = 0.006350235899329861
$ws.Range("A4").Value = "year"
$ws.Range("B4").Value = 0.0004615048445069547
$ws.Range("C4").Value = 0.001409832983329938
$ws.Range("A5").Value = "COVID19"
$ws.Range("B5").Value = -0.00005123672521341049
$ws.Range("C5").Value = 0.00001028609278392733
$ws.Range("A6").Value = "AR2"
$ws.Range("B6").Value = -0.001257121757073776
$ws.Range("C6").Value = 0.003761903464868745
$ws.Range("A7").Value = "monthcos"
$ws.Range("B7").Value = -0.002183523881799499
$ws.Range("C7").Value = 0.002938459148310781
$ws.Range("A8").Value = "AR12"
$ws.Range("B8").Value = -0.02008014542261183
$ws.Range("C8").Value = 0.01926523909363137
$ws.Range("A9").Value = "monthsin"
$ws.Range("B9").Value = -0.04703717703343448
$ws.Range("C9").Value = 0.003670901003530233
$ws.Range("A10").Value = "AR3"
$ws.Range("B10").Value = -0.05731880094367767
$ws.Range("C10").Value = 0.007064769973646699
$ws.Range("A11").Value = "AR1"
$ws.Range("B11").Value = -0.1258926100592255
$ws.Range("C11").Value = 0.008154204759094777
$ws.Range("A12").Value = "AR24"
$ws.Range("B12").Value = -0.1261745027946785
$ws.Range("C12").Value = 0.03815650035522476
$ws.Range("A13").Value = "AR4"
$ws.Range("B13").Value = -0.1333134038443707
$ws.Range("C13").Value = 0.02440575590958999

# --- mlp_feature_importance ---
$ws = $wb.Worksheets.Item("mlp_feature_importance")
$ws.Range("A2").Value = "AR24"
$ws.Range("B2").Value = 0.007250122454850549
$ws.Range("C2").Value = 0.001670788838701229
$ws.Range("A3").Value = "AR12"
$ws.Range("B3").Value = 0.006233948449155102
$ws.Range("C3").Value = 0.0008851076106442664
$ws.Range("A4").Value = "AR3"
$ws.Range("B4").Value = 0.003906327366523232
$ws.Range("C4").Value = 0.001497029477093239
$ws.Range("A5").Value = "AR4"
$ws.Range("B5").Value = 0.003619988228748961
$ws.Range("C5").Value = 0.001019337869048467
$ws.Range("A6").Value = "AR2"
$ws.Range("B6").Value = 0.0006562745151698745
$ws.Range("C6").Value = 0.0004641032433140799
$ws.Range("A7").Value = "AR1"
$ws.Range("B7").Value = 0.0004567229677268081
$ws.Range("C7").Value = 0.0007514655170052605
$ws.Range("A8").Value = "monthcos"
$ws.Range("B8").Value = 0.0001503576099719073
$ws.Range("C8").Value = 0.0001124937357048963
$ws.Range("A9").Value = "monthsin"
$ws.Range("B9").Value = 0.00004157687954298162
$ws.Range("C9").Value = 0.00005448701086537222
$ws.Range("A10").Value = "COVID19"
$ws.Range("B10").Value = 0.00001653001330601533
$ws.Range("C10").Value = 0.00001098612357695792
$ws.Range("A11").Value = "year"
$ws.Range("B11").Value = -0.0000303191404348091
$ws.Range("C11").Value = 0.0002869721286932487
$ws.Range("A12").Value = "t_COVID19"
$ws.Range("B12").Value = -0.0005596661944410553
$ws.Range("C12").Value = 0.001444589027826579
$ws.Range("A13").Value = "t"
$ws.Range("B13").Value = -0.01728425285675226
$ws.Range("C13").Value = 0.004994381227968889

# --- silverkite_summary_stats ---
$ws = $wb.Worksheets.Item("silverkite_summary_stats")
$ws.Range("C2").Value = 1.544421717023353
$ws.Range("D2").Value = 0.094
$ws.Range("F2").Value = "[-0.4564361494213066, 5.792062201414955]"
$ws.Range("C3").Value = 0.3891732039298657
$ws.Range("D3").Value = 0.992
$ws.Range("F3").Value = "[-0.7753360262476652, 0.7627512190636426]"
$ws.Range("C4").Value = 0.5425205455871679
$ws.Range("D4").Value = 0.708
$ws.Range("F4").Value = "[-0.755611253327032, 1.2861015085581704]"
$ws.Range("C5").Value = 0.5899332314831082
$ws.Range("D5").Value = 0.334
$ws.Range("F5").Value = "[-1.618401886487736, 0.6883146855985145]"
$ws.Range("C6").Value = 0.4191992903726603
$ws.Range("D6").Value = 0.716
$ws.Range("F6").Value = "[-0.9449003024563061, 0.6749149016616901]"
$ws.Range("C7").Value = 0.4288897016400642
$ws.Range("D7").Value = 0.804
$ws.Range("F7").Value = "[-1.0024380272678148, 0.7217601494855292]"
$ws.Range("C13").Value = 0.5511932688032536
$ws.Range("D13").Value = 0.116
$ws.Range("F13").Value = "[-0.020132685133808816, 1.8643070772777819]"
$ws.Range("C14").Value = 0.2332180076451266
$ws.Range("D14").Value = 0.986
$ws.Range("F14").Value = "[-0.4214042642518702, 0.5168253223655759]"
$ws.Range("C15").Value = 0.815659438792821
$ws.Range("D15").Value = 0.148
$ws.Range("F15").Value = "[0.0, 2.5496134031066644]"
$ws.Range("C16").Value = 0.2020242028099362
$ws.Range("D16").Value = 0.1
$ws.Range("F16").Value = "[0.0, 0.7171409344397852]"
$ws.Range("C17").Value = 0.4196260890345305
$ws.Range("D17").Value = 0.152
$ws.Range("F17").Value = "[0.0, 1.3886379193029224]"
$ws.Range("C23").Value = 0.9980553945503653
$ws.Range("D23").Value = 0.0
$ws.Range("F23").Value = "[1.462477768980686, 5.262041347702955]"
$ws.Range("C24").Value = 0.6770099045916615
$ws.Range("D24").Value = 0.008
$ws.Range("E24").Value = "**"
$ws.Range("F24").Value = "[-3.007692376205913, -0.4181730986133249]"
$ws.Range("C25").Value = 1.325669925079384
$ws.Range("D25").Value = 0.008
$ws.Range("E25").Value = "**"
$ws.Range("F25").Value = "[-5.7649906126828245, -0.6337473432094788]"
$ws.Range("C31").Value = 0.6130857228646857
$ws.Range("D31").Value = 0.382
$ws.Range("F31").Value = "[-0.5854637770020035, 1.8580668100690532]"
$ws.Range("C32").Value = 0.7563246350007
$ws.Range("D32").Value = 0.018
$ws.Range("F32").Value = "[-3.197839672329523, -0.4792585720084501]"
$ws.Range("C33").Value = 0.9980411996827439
$ws.Range("D33").Value = 0.92
$ws.Range("F33").Value = "[-1.9664996320630035, 2.008529718296263]"
$ws.Range("C38").Value = 0.9248778903585332
$ws.Range("D38").Value = 0.0
$ws.Range("F38").Value = "[-5.226158040232855, -1.678895748028335]"
$ws.Range("C39").Value = 1.024863760466941
$ws.Range("D39").Value = 0.0
$ws.Range("F39").Value = "[-8.605459576591715, -4.593459391260007]"
$ws.Range("C40").Value = 1.041776975603151
$ws.Range("D40").Value = 0.0
$ws.Range("F40").Value = "[2.7573075034947863, 6.725820003852739]"
$ws.Range("C41").Value = 1.165834823893171
$ws.Range("D41").Value = 0.0
$ws.Range("F41").Value = "[-9.49152234496447, -4.828637851555698]"
$ws.Range("C42").Value = 1.114592134297336
$ws.Range("D42").Value = 0.108
$ws.Range("F42").Value = "[-0.3070192267974259, 3.863262132280031]"
$ws.Range("C53").Value = 0.03254632962303766
$ws.Range("D53").Value = 0.658
$ws.Range("F53").Value = "[-0.05173832225043615, 0.08094340731457195]"
